# Atualizacao de bases das ligas (Germany Landesliga), do dia: 13-06-2024 as 19:35
#
# The underlying source refresh re-pulled match rows 11-13 and 83-86 in a
# different order (the match ids/teams/odds are unchanged as data, only
# their row position shuffled), so the fix is to rewrite B:AD for those
# seven rows with their new row assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($ws, $rowNum, $values) {
    $arr = New-Object "object[,]" 1,$values.Count
    for ($i = 0; $i -lt $values.Count; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $startCol = 2  # column B
    $endCol = $startCol + $values.Count - 1
    $rng = $ws.Range($ws.Cells.Item($rowNum, $startCol), $ws.Cells.Item($rowNum, $endCol))
    $rng.Value = $arr
}

Set-RowValues $ws 11 @(7035048, "Germany Landesliga", 45147.625, "SG Unterrath", "TuRU Dsseldorf", 1, 0, 1, 0, "H", 3.25, 4, 1.8, 2.9, 4, 1.95, 0.5, 1.8, 2, 3, 1.75, 1.95, 1.9, -1, -1, 0.8, -1, -1, 0.95)
Set-RowValues $ws 12 @(7035046, "Germany Landesliga", 45147.625, "Cronenberger SC", "FC Viersen", 0, 2, 0, 1, "A", 2, 3.6, 3, 2, 3.6, 3, -0.25, 1.8, 2, 2.75, 1.8, 2, -1, -1, 2, -1, 1, -1, 1)
Set-RowValues $ws 13 @(7035047, "Germany Landesliga", 45147.625, "SC Dsseldorf West", "VfL Viktoria JuchenGarzweiler", 3, 4, 3, 1, "A", 1.909, 3.75, 3.1, 2.2, 3.6, 2.625, -0.25, 2, 1.8, 3, 1.825, 1.975, -1, -1, 1.625, -1, 0.8, 0.825, -1)
Set-RowValues $ws 83 @(8075530, "Germany Landesliga", 45392.60416666666, "TuRU Dsseldorf", "DV Solingen", 2, 0, 1, 0, "H", 2.1, 3.75, 2.7, 2.375, 3.75, 2.45, 0, 1.85, 1.95, 3, 1.85, 1.95, 1.375, -1, -1, 0.8500000000000001, -1, -1, 0.95)
Set-RowValues $ws 84 @(8075296, "Germany Landesliga", 45392.60416666666, "FC Monheim", "VFB Hilden II", 1, 2, 1, 2, "A", 1.533, 4.75, 4, 1.4, 5.25, 5, -1.5, 1.975, 1.825, 3.75, 1.9, 1.9, -1, -1, 4, -1, 0.825, -1, 0.8999999999999999)
Set-RowValues $ws 85 @(8075593, "Germany Landesliga", 45392.625, "FSV Duisburg", "ESC Rellinghausen", 0, 2, 0, 1, "A", 6.5, 4.5, 1.363, 4.75, 4.5, 1.5, 1.25, 1.875, 1.925, 3.75, 2, 1.8, -1, -1, 0.5, -1, 0.925, -1, 0.8)
Set-RowValues $ws 86 @(8075670, "Germany Landesliga", 45392.625, "SF Niederwenigern", "VfB Frohnhausen", 2, 1, 0, 1, "H", 1.222, 6.5, 9.5, 1.181, 7, 11, -2.5, 1.925, 1.875, 4.25, 1.8, 2, 0.181, -1, -1, -1, 0.875, -1, 1)
